# Automatic update of files.
#
# Applies the source-system refresh to rows 2-5 and 7-9 of the Artfynd
# sightings export: new report Ids, corrected locality name, tightened
# accuracy radius, shifted dates, dropped public comment, and a change of
# observer/reporter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper so date-like text (e.g. "2017-08-30") is written back as literal
# text instead of being auto-coerced into an Excel date serial number.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

$rows = 2,3,4,5,7,8,9

foreach ($r in $rows) {
    $ws.Range("P$r").Value = "Njuöniesvarie, Ly lm"
    $ws.Range("S$r").Value = 10
    Set-TextValue $ws.Range("Y$r") "2017-08-30"
    Set-TextValue $ws.Range("AA$r") "2017-08-30"
    $ws.Range("AW$r").Value = "Jonas Nordenström"
    $ws.Range("AX$r").Value = "Jonas Nordenström"
}

# New report Ids (column A)
$ws.Range("A2").Value = 79087787
$ws.Range("A3").Value = 79087818
$ws.Range("A4").Value = 79087802
$ws.Range("A5").Value = 79087805
$ws.Range("A7").Value = 79087826
$ws.Range("A8").Value = 79087808
$ws.Range("A9").Value = 79087821

# Newly-reported activity for row 2
$ws.Range("M2").Value = "färska spår"

# Public comment ("ringade träd") withdrawn on rows 2 and 9
$ws.Range("AC2").Value = ""
$ws.Range("AC9").Value = ""
